$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (never auto-converted to a number)
# by building it via a temporary formula cell and pasting VALUES ONLY,
# which avoids Excel's "smart" numeric-string coercion on Range.Value
# while leaving cell formatting/style completely untouched.
function Set-TextValue($cellRef, [string]$text) {
    $ws.Range("ZZ1").Formula = '="' + $text + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("D2").Value = '48.168.73'
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").Value = '2.520.90'
$ws.Range("E3").Value = '  +1.29%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue "D5" "110.09"
$ws.Range("E5").Value = '  +1.28%  '

Set-TextValue "D6" "323.07"
$ws.Range("E6").Value = '  +0.35%  '

Set-TextValue "D7" "0.533"
$ws.Range("E7").Value = '  +2.02%  '

$ws.Range("E8").Value = '  +0.07%  '

Set-TextValue "D9" "0.555"
$ws.Range("E9").Value = '  +4.06%  '

Set-TextValue "D10" "40.69"
$ws.Range("E10").Value = '  +4.86%  '

Set-TextValue "D11" "20.48"
$ws.Range("E11").Value = '  +12.46%  '

Set-TextValue "D12" "0.0827"
$ws.Range("E12").Value = '  +2.14%  '

$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("E14").Value = '  +2.02%  '

$ws.Range("D15").Value = '2.917.86'
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").Value = '2.519.84'
$ws.Range("E16").Value = '  +1.25%  '

Set-TextValue "D17" "0.856"
$ws.Range("E17").Value = '  +1.30%  '

$ws.Range("D18").Value = '48.015.58'
$ws.Range("E18").Value = '  +2.08%  '

$ws.Range("E19").Value = '  +3.88%  '

Set-TextValue "D20" "6.63"
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("D21").Value = '0.0₃0949'
$ws.Range("E21").Value = '  +1.51%  '

Set-TextValue "D22" "2.70"
$ws.Range("E22").Value = '  -2.06%  '

Set-TextValue "D23" "72.07"
$ws.Range("E23").Value = '  +2.09%  '

Set-TextValue "D24" "263.83"
$ws.Range("E24").Value = '  +7.33%  '

$ws.Range("E25").Value = '  +0.39%  '

Set-TextValue "D26" "26.16"
$ws.Range("E26").Value = '  +1.62%  '

$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("E29").Value = '  +3.93%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D30" "36.60"
$ws.Range("E30").Value = '  +4.47%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D31" "2.20"
$ws.Range("E31").Value = '  -3.61%  '

Set-TextValue "D32" "49.73"
$ws.Range("E32").Value = '  -0.36%  '

Set-TextValue "D33" "19.97"
$ws.Range("E33").Value = '  +0.00%  '

Set-TextValue "D34" "5.41"
$ws.Range("E34").Value = '  +0.30%  '

$ws.Range("E35").Value = '  +0.02%  '

Set-TextValue "D36" "0.0792"
$ws.Range("E36").Value = '  +1.25%  '

$ws.Range("E37").Value = '  +1.47%  '

$ws.Range("E38").Value = '  +1.48%  '

$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("E40").Value = '  +0.80%  '

Set-TextValue "D41" "120.53"
$ws.Range("E41").Value = '  +1.42%  '

$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D42" "2.20"
$ws.Range("E42").Value = '  -0.87%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D43" "21.93"
$ws.Range("E43").Value = '  +3.26%  '

$ws.Range("D45").Value = '2.022.37'
$ws.Range("E45").Value = '  +2.16%  '

$ws.Range("E46").Value = '  +4.95%  '

Set-TextValue "D47" "1.92"
$ws.Range("E47").Value = '  +8.65%  '

$ws.Range("E48").Value = '  +1.17%  '

Set-TextValue "D49" "9.13"
$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("E50").Value = '  +2.85%  '

Set-TextValue "D51" "78.90"
$ws.Range("E51").Value = '  +2.53%  '

$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = 0
